$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.35379"
$ws.Range("H2").Value = [double]"1.06137"
$ws.Range("I2").Value = [double]"0.008539219707589354"
$ws.Range("J2").Value = [double]"0.008539219707589352"
$ws.Range("M2").Value = [double]"7.579746333333333"
$ws.Range("N2").Value = [double]"22.739239"
$ws.Range("O2").Value = [double]"0.0686314777863378"
$ws.Range("P2").Value = [double]"0.0686314777863378"
$ws.Range("Q2").Value = [double]"2.68163845527"
$ws.Range("R2").Value = [double]"24.13474609743"
$ws.Range("S2").Value = [double]"0.0005860592676740767"
$ws.Range("T2").Value = [double]"0.0005860592676740766"

$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.35379"
$ws.Range("H3").Value = [double]"1.06137"
$ws.Range("I3").Value = [double]"0.008539219707589354"
$ws.Range("J3").Value = [double]"0.008539219707589352"
$ws.Range("O3").Value = [double]"0.0596740760116217"
$ws.Range("P3").Value = [double]"0.05967407601162171"
$ws.Range("Q3").Value = [double]"2.331645801270001"
$ws.Range("R3").Value = [double]"20.98481221143"
$ws.Range("S3").Value = [double]"0.0005095700459106251"
$ws.Range("T3").Value = [double]"0.0005095700459106251"

$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.35379"
$ws.Range("H4").Value = [double]"1.06137"
$ws.Range("I4").Value = [double]"0.008539219707589354"
$ws.Range("J4").Value = [double]"0.008539219707589352"
$ws.Range("M4").Value = [double]"96.08192699999999"
$ws.Range("N4").Value = [double]"288.245781"
$ws.Range("O4").Value = [double]"0.8699822327258658"
$ws.Range("P4").Value = [double]"0.8699822327258659"
$ws.Range("Q4").Value = [double]"33.99282495333"
$ws.Range("R4").Value = [double]"305.93542457997"
$ws.Range("S4").Value = [double]"0.007428969426945302"
$ws.Range("T4").Value = [double]"0.007428969426945301"

$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"0.6666666666666666"
$ws.Range("G5").Value = [double]"0.35379"
$ws.Range("H5").Value = [double]"1.06137"
$ws.Range("I5").Value = [double]"0.008539219707589354"
$ws.Range("J5").Value = [double]"0.008539219707589352"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.189099"
$ws.Range("N5").Value = [double]"0.5672970000000001"
$ws.Range("O5").Value = [double]"0.001712213476174646"
$ws.Range("P5").Value = [double]"0.001712213476174646"
$ws.Range("Q5").Value = [double]"0.06690133521000001"
$ws.Range("R5").Value = [double]"0.6021120168900002"
$ws.Range("S5").Value = [double]"1.462096705935061E-05"
$ws.Range("T5").Value = [double]"1.462096705935061E-05"

$ws.Range("G6").Value = [double]"36.97491766666666"
$ws.Range("I6").Value = [double]"0.8924416903408624"
$ws.Range("J6").Value = [double]"0.8924416903408623"
$ws.Range("M6").Value = [double]"7.579746333333333"
$ws.Range("N6").Value = [double]"22.739239"
$ws.Range("O6").Value = [double]"0.0686314777863378"
$ws.Range("P6").Value = [double]"0.0686314777863378"
$ws.Range("Q6").Value = [double]"280.2604966092185"
$ws.Range("R6").Value = [double]"2522.344469482967"
$ws.Range("S6").Value = [double]"0.06124959204623065"
$ws.Range("T6").Value = [double]"0.06124959204623064"

$ws.Range("G7").Value = [double]"36.97491766666666"
$ws.Range("I7").Value = [double]"0.8924416903408624"
$ws.Range("J7").Value = [double]"0.8924416903408623"
$ws.Range("O7").Value = [double]"0.0596740760116217"
$ws.Range("P7").Value = [double]"0.05967407601162171"
$ws.Range("S7").Value = [double]"0.05325563326534078"
$ws.Range("T7").Value = [double]"0.05325563326534078"

$ws.Range("G8").Value = [double]"36.97491766666666"
$ws.Range("I8").Value = [double]"0.8924416903408624"
$ws.Range("J8").Value = [double]"0.8924416903408623"
$ws.Range("M8").Value = [double]"96.08192699999999"
$ws.Range("N8").Value = [double]"288.245781"
$ws.Range("O8").Value = [double]"0.8699822327258658"
$ws.Range("P8").Value = [double]"0.8699822327258659"
$ws.Range("Q8").Value = [double]"3552.621340079676"
$ws.Range("R8").Value = [double]"31973.59206071709"
$ws.Range("S8").Value = [double]"0.7764084143403892"
$ws.Range("T8").Value = [double]"0.7764084143403892"

$ws.Range("G9").Value = [double]"36.97491766666666"
$ws.Range("I9").Value = [double]"0.8924416903408624"
$ws.Range("J9").Value = [double]"0.8924416903408623"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.189099"
$ws.Range("N9").Value = [double]"0.5672970000000001"
$ws.Range("O9").Value = [double]"0.001712213476174646"
$ws.Range("P9").Value = [double]"0.001712213476174646"
$ws.Range("Q9").Value = [double]"6.991919955849"
$ws.Range("R9").Value = [double]"62.927279602641"
$ws.Range("S9").Value = [double]"0.001528050688901705"
$ws.Range("T9").Value = [double]"0.001528050688901705"

$ws.Range("G10").Value = [double]"4.102478333333333"
$ws.Range("H10").Value = [double]"12.307435"
$ws.Range("I10").Value = [double]"0.09901908995154843"
$ws.Range("J10").Value = [double]"0.0990190899515484"
$ws.Range("M10").Value = [double]"7.579746333333333"
$ws.Range("N10").Value = [double]"22.739239"
$ws.Range("O10").Value = [double]"0.0686314777863378"
$ws.Range("P10").Value = [double]"0.0686314777863378"
$ws.Range("Q10").Value = [double]"31.09574510466278"
$ws.Range("R10").Value = [double]"279.861705941965"
$ws.Range("S10").Value = [double]"0.00679582647243308"
$ws.Range("T10").Value = [double]"0.006795826472433078"

$ws.Range("G11").Value = [double]"4.102478333333333"
$ws.Range("H11").Value = [double]"12.307435"
$ws.Range("I11").Value = [double]"0.09901908995154843"
$ws.Range("J11").Value = [double]"0.0990190899515484"
$ws.Range("O11").Value = [double]"0.0596740760116217"
$ws.Range("P11").Value = [double]"0.05967407601162171"
$ws.Range("Q11").Value = [double]"27.03730003877389"
$ws.Range("R11").Value = [double]"243.335700348965"
$ws.Range("S11").Value = [double]"0.005908872700370308"
$ws.Range("T11").Value = [double]"0.005908872700370307"

$ws.Range("G12").Value = [double]"4.102478333333333"
$ws.Range("H12").Value = [double]"12.307435"
$ws.Range("I12").Value = [double]"0.09901908995154843"
$ws.Range("J12").Value = [double]"0.0990190899515484"
$ws.Range("M12").Value = [double]"96.08192699999999"
$ws.Range("N12").Value = [double]"288.245781"
$ws.Range("O12").Value = [double]"0.8699822327258658"
$ws.Range("P12").Value = [double]"0.8699822327258659"
$ws.Range("Q12").Value = [double]"394.1740237424149"
$ws.Range("R12").Value = [double]"3547.566213681734"
$ws.Range("S12").Value = [double]"0.08614484895853144"
$ws.Range("T12").Value = [double]"0.08614484895853143"

$ws.Range("G13").Value = [double]"4.102478333333333"
$ws.Range("H13").Value = [double]"12.307435"
$ws.Range("I13").Value = [double]"0.09901908995154843"
$ws.Range("J13").Value = [double]"0.0990190899515484"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.189099"
$ws.Range("N13").Value = [double]"0.5672970000000001"
$ws.Range("O13").Value = [double]"0.001712213476174646"
$ws.Range("P13").Value = [double]"0.001712213476174646"
$ws.Range("Q13").Value = [double]"0.775774550355"
$ws.Range("R13").Value = [double]"6.981970953195"
$ws.Range("S13").Value = [double]"0.0001695418202135907"
$ws.Range("T13").Value = [double]"0.0001695418202135906"
